$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing date entries (rows 8-10) to correct November dates
$ws.Range("A8").Value = 41956
$ws.Range("A9").Value = 41957
$ws.Range("A10").Value = 41957

# Add new row 11 for work done on 22/11/2014
$ws.Range("A11").Value = 41965
$ws.Range("B11").Value = 2.5
$ws.Range("C11").Value = "Footer navigation and mobile responsive behavior for footer."

# Update the selected cell shown when the sheet is opened
[void]$ws.Range("C12").Select()
